# Update the notes (cell comments) to the column keys, adjust the header
# cell style for C1:D1 (left-aligned now), set the best-fit-ish column
# widths for A:H, and move the active selection to G2 - as captured by the
# commit "Update notes to column keys in metadata files (Excel)".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

# --- A1 comment -----------------------------------------------------------
# Merge the former B1 note into A1 and reword it.
$a1Comment = $ws.Range("A1").Comment
$a1Text = "Unique sample IDs " + $nl + "(Enter ""?setup_expts"" under an R console for details of the column keys.)"
$a1Comment.Text($a1Text)

# --- B1 comment -------------------------------------------------------------
# The separate B1 note is removed now that its content lives on A1.
$ws.Range("B1").Comment.Delete() | Out-Null

# --- E1 comment -------------------------------------------------------------
$e1Comment = $ws.Range("E1").Comment
$e1Text = "The file names of RAW MS data" + $nl + "(Enter alternatively the file names in frac_smry.xlsx if more than one RAW file per TMT set)"
$e1Comment.Text($e1Text)

# --- F1 comment -------------------------------------------------------------
$f1Comment = $ws.Range("F1").Comment
$f1Comment.Text("Non-void character strings to indicate refernce channels.")

# --- G1 comment -------------------------------------------------------------
$g1Comment = $ws.Range("G1").Comment
$g1Text = "Samples to be selected for analyses " + $nl + "(Enter ""?setup_expts"" under an R console for details of the column keys.)"
$g1Comment.Text($g1Text)

# --- Header style for C1:D1 --------------------------------------------------
# Left-align these header cells (keeps their existing fill/border/font),
# which Excel records as a new cellXfs entry.
$ws.Range("C1:D1").HorizontalAlignment = -4131 ## xlLeft

# --- Column widths A:H --------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 18.666666666666668
$ws.Columns("B:B").ColumnWidth = 8.166666666666666
$ws.Columns("C:C").ColumnWidth = 7.666666666666667
$ws.Columns("D:D").ColumnWidth = 8.166666666666666
$ws.Columns("E:E").ColumnWidth = 7.666666666666667
$ws.Columns("F:F").ColumnWidth = 8.166666666666666
$ws.Columns("G:G").ColumnWidth = 18.666666666666668
$ws.Columns("H:H").ColumnWidth = 14.5

# --- Selection ---------------------------------------------------------------
$ws.Range("G2").Select() | Out-Null
